$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.41746662583399
$ws.Range("C2").Value = 16.72799125855175
$ws.Range("D2").Value = 6.596485066606035
$ws.Range("E2").Value = 12.54334689024369
$ws.Range("F2").Value = 45.95126704239517
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 29.1645158382879
$ws.Range("J2").Value = 10.34759226140388

$ws.Range("B3").Value = 18.95060772146521
$ws.Range("C3").Value = 16.24472651386505
$ws.Range("D3").Value = 6.599381963312862
$ws.Range("E3").Value = 12.52409565201733
$ws.Range("F3").Value = 45.68090976712895
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 29.12847201596303
$ws.Range("J3").Value = 10.35653013457559

$ws.Range("B4").Value = 18.66409336860236
$ws.Range("C4").Value = 15.94628212193704
$ws.Range("D4").Value = 6.601847285039836
$ws.Range("E4").Value = 12.51497466854458
$ws.Range("F4").Value = 45.52828976688014
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 29.11418353287728
$ws.Range("J4").Value = 10.36396920002342

$ws.Range("B5").Value = 18.54756221238148
$ws.Range("C5").Value = 15.82444127546175
$ws.Range("D5").Value = 6.603025774198426
$ws.Range("E5").Value = 12.51193858857168
$ws.Range("F5").Value = 45.46949389641871
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 29.11032478413129
$ws.Range("J5").Value = 10.36749034192591

$ws.Range("B6").Value = 18.52823162319701
$ws.Range("C6").Value = 15.80420259586319
$ws.Range("D6").Value = 6.603231996308604
$ws.Range("E6").Value = 12.51147561726947
$ws.Range("F6").Value = 45.45993697096493
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 29.10980238550898
$ws.Range("J6").Value = 10.36810456805997

$ws.Range("B7").Value = 18.66252061825097
$ws.Range("C7").Value = 15.9446395492065
$ws.Range("D7").Value = 6.60186247306553
$ws.Range("E7").Value = 12.51493096408207
$ws.Range("F7").Value = 45.52748302721321
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 29.11412355259153
$ws.Range("J7").Value = 10.36401470612511

$ws.Range("B8").Value = 19.25657748970477
$ws.Range("C8").Value = 16.56184055656843
$ws.Range("D8").Value = 6.597342129388995
$ws.Range("E8").Value = 12.53615017023911
$ws.Range("F8").Value = 45.85529265750294
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 29.15045581235246
$ws.Range("J8").Value = 10.35026852355536

$ws.Range("B9").Value = 20.41457352546956
$ws.Range("C9").Value = 17.74990608433522
$ws.Range("D9").Value = 6.593873206090761
$ws.Range("E9").Value = 12.5990701093564
$ws.Range("F9").Value = 46.6025389063586
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 29.28421873171781
$ws.Range("J9").Value = 10.33883335615768

$ws.Range("B10").Value = 21.25092261244774
$ws.Range("C10").Value = 18.59854820119009
$ws.Range("D10").Value = 6.594543222289396
$ws.Range("E10").Value = 12.65811267732865
$ws.Range("F10").Value = 47.21245118728808
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 29.42088862709652
$ws.Range("J10").Value = 10.33994146324607

$ws.Range("B11").Value = 21.62630319079016
$ws.Range("C11").Value = 18.97737736105558
$ws.Range("D11").Value = 6.595531666670415
$ws.Range("E11").Value = 12.68770771906756
$ws.Range("F11").Value = 47.50244644421336
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 29.49142252029701
$ws.Range("J11").Value = 10.34251771449943

$ws.Range("B12").Value = 21.76757087484813
$ws.Range("C12").Value = 19.119645795421
$ws.Range("D12").Value = 6.596002961247622
$ws.Range("E12").Value = 12.69930318441807
$ws.Range("F12").Value = 47.61399790117384
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 29.5193328616526
$ws.Range("J12").Value = 10.34379153494249

$ws.Range("B13").Value = 21.7371878002143
$ws.Range("C13").Value = 19.08906068610657
$ws.Range("D13").Value = 6.595897166878228
$ws.Range("E13").Value = 12.69678870121604
$ws.Range("F13").Value = 47.58989717758563
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 29.51326851754302
$ws.Range("J13").Value = 10.34350392750893

$ws.Range("B14").Value = 21.63794384215976
$ws.Range("C14").Value = 18.98910643339971
$ws.Range("D14").Value = 6.595568504078248
$ws.Range("E14").Value = 12.68865392994188
$ws.Range("F14").Value = 47.51158944575306
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 29.49369467417633
$ws.Range("J14").Value = 10.3426165340985

$ws.Range("B15").Value = 21.57703507463602
$ws.Range("C15").Value = 18.9277229870329
$ws.Range("D15").Value = 6.595379780572258
$ws.Range("E15").Value = 12.68372158049085
$ws.Range("F15").Value = 47.46384778618598
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 29.4818614100503
$ws.Range("J15").Value = 10.34211182637187

$ws.Range("B16").Value = 21.22627528954488
$ws.Range("C16").Value = 18.57363257334837
$ws.Range("D16").Value = 6.594492269745612
$ws.Range("E16").Value = 12.65623321557189
$ws.Range("F16").Value = 47.19374607072802
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 29.41644722949617
$ws.Range("J16").Value = 10.33981482364702

$ws.Range("B17").Value = 21.00968407552715
$ws.Range("C17").Value = 18.35445042443675
$ws.Range("D17").Value = 6.594121978002433
$ws.Range("E17").Value = 12.64006741122453
$ws.Range("F17").Value = 47.03121358194267
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 29.37845881784987
$ws.Range("J17").Value = 10.33893666528562

$ws.Range("B18").Value = 20.88463678269206
$ws.Range("C18").Value = 18.22771075348651
$ws.Range("D18").Value = 6.593973424025161
$ws.Range("E18").Value = 12.63102715083632
$ws.Range("F18").Value = 46.9389141112319
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 29.3573961932359
$ws.Range("J18").Value = 10.33862660786143

$ws.Range("B19").Value = 20.84222200686634
$ws.Range("C19").Value = 18.18468820046666
$ws.Range("D19").Value = 6.593934231644681
$ws.Range("E19").Value = 12.62801071029356
$ws.Range("F19").Value = 46.90786858371656
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 29.35039999359737
$ws.Range("J19").Value = 10.33855511434695

$ws.Range("B20").Value = 21.03279032380277
$ws.Range("C20").Value = 18.37785338911557
$ws.Range("D20").Value = 6.594154737623609
$ws.Range("E20").Value = 12.64176163357368
$ws.Range("F20").Value = 47.04839328150271
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 29.38242127288035
$ws.Range("J20").Value = 10.33900995875591

$ws.Range("B21").Value = 21.66711924577202
$ws.Range("C21").Value = 19.01849871185975
$ws.Range("D21").Value = 6.595662418526989
$ws.Range("E21").Value = 12.69103280890239
$ws.Range("F21").Value = 47.53454375060902
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 29.49941142358735
$ws.Range("J21").Value = 10.34286908724339

$ws.Range("B22").Value = 22.07649414427978
$ws.Range("C22").Value = 19.43022626763608
$ws.Range("D22").Value = 6.597212499407651
$ws.Range("E22").Value = 12.72549598711753
$ws.Range("F22").Value = 47.86235921231881
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 29.58286762590106
$ws.Range("J22").Value = 10.34712969354214

$ws.Range("B23").Value = 21.85852524198411
$ws.Range("C23").Value = 19.21116248563492
$ws.Range("D23").Value = 6.596333958065245
$ws.Range("E23").Value = 12.70689717935936
$ws.Range("F23").Value = 47.68649808026261
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 29.53768638124903
$ws.Range("J23").Value = 10.34469661102914

$ws.Range("B24").Value = 21.02234561615958
$ws.Range("C24").Value = 18.36727516791673
$ws.Range("D24").Value = 6.594139726520769
$ws.Range("E24").Value = 12.64099488513186
$ws.Range("F24").Value = 47.0406227752761
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 29.38062742536151
$ws.Range("J24").Value = 10.33897621601867

$ws.Range("B25").Value = 20.10313771746582
$ws.Range("C25").Value = 17.43203923071045
$ws.Range("D25").Value = 6.594241253644697
$ws.Range("E25").Value = 12.5797833842682
$ws.Range("F25").Value = 46.38947200627952
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 29.24130302879481
$ws.Range("J25").Value = 10.34025974244614
